$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Python")
$ws2 = $wb.Worksheets.Item("SQL")

# --- SQL sheet: D49 gets a literal 0 (style/number format unchanged) ---
$ws2.Range("D49").Value = 0

# --- Python sheet: mark newly-completed topics (COMPLETED + date), rows grouped by day ---
$ws1.Range("B2:C2").Copy()
$ws1.Range("B13:C21").PasteSpecial(-4122)
$ws1.Range("B13:B21").Value = "COMPLETED"
$ws1.Range("C13:C21").Value = 45870
$ws1.Range("B2:C2").Copy()
$ws1.Range("B23:C32").PasteSpecial(-4122)
$ws1.Range("B23:B32").Value = "COMPLETED"
$ws1.Range("C23:C32").Value = 45871
$ws1.Range("B2:C2").Copy()
$ws1.Range("B34:C43").PasteSpecial(-4122)
$ws1.Range("B34:B43").Value = "COMPLETED"
$ws1.Range("C34:C43").Value = 45871
$ws1.Range("B2:C2").Copy()
$ws1.Range("B45:C54").PasteSpecial(-4122)
$ws1.Range("B45:B54").Value = "COMPLETED"
$ws1.Range("C45:C54").Value = 45872
$ws1.Range("B2:C2").Copy()
$ws1.Range("B56:C65").PasteSpecial(-4122)
$ws1.Range("B56:B65").Value = "COMPLETED"
$ws1.Range("C56:C65").Value = 45873
$ws1.Range("B2:C2").Copy()
$ws1.Range("B67:C74").PasteSpecial(-4122)
$ws1.Range("B67:B74").Value = "COMPLETED"
$ws1.Range("C67:C74").Value = 45874
$ws1.Range("B2:C2").Copy()
$ws1.Range("B78:C82").PasteSpecial(-4122)
$ws1.Range("B78:B82").Value = "COMPLETED"
$ws1.Range("C78:C82").Value = 45875

$excel.CutCopyMode = $false

# --- Column A on Python sheet was widened (best-fit) ---
$ws1.Columns.Item(1).ColumnWidth = 51

# --- Window/view state: SQL tab scrolled/zoomed/selected first, then focus returned to Python tab ---
$ws2.Activate()
$excel.ActiveWindow.Zoom = 88
$ws2.Range("D49").Select()

$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.Zoom = 98
$ws1.Range("D80").Select()
